$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as plain text values,
# matching the original inline-string cell formatting.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.408.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.819.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.08%  '
$ws.Range('E4').Value = '  -1.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '332.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('E6').Value = '  -0.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4552'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3823'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.88'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9557'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.91'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.38%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.816'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.17%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.808.02'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.029'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06579'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001016'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.398.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.256'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.77'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.259'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.59'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.018.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.22'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.033'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.251'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.55'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09299'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9245'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.567'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.202'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.312'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05890'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02179'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.29%  '
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.037'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.133'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5703'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1809'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.854'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.262'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5370'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.66'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.868'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06840'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '109.75'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.77%  '
$ws.Range('E51').Value = '  -33.12%  '
